$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# "b.md" has finished its handback cycle and is now ready to be handed off
# again, so its status flips from "Handed back: in sync with en-US" to
# "Ready for handoff" on every sheet, and the zh-cn / de-de detail sheets
# get a new handoff file + timestamp recorded for that row.
# ---------------------------------------------------------------------------

# --- Overview sheet: update the b.md row (row 3) status in both language columns
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = "Ready for handoff"
$wsOverview.Range("C3").Value = "Ready for handoff"

# --- zh-cn sheet: update b.md row (row 3)
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("B3").Value = "Ready for handoff"
$wsZh.Range("C3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$wsZh.Range("D3").Value = "2016-03-08 08:21:17"

foreach ($hl in $wsZh.Hyperlinks) {
    if ($hl.Range.Address() -eq '$C$3') {
        $hl.TextToDisplay = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
    }
}

# --- de-de sheet: update b.md row (row 3)
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("B3").Value = "Ready for handoff"
$wsDe.Range("C3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$wsDe.Range("D3").Value = "2016-03-08 08:21:20"

foreach ($hl in $wsDe.Hyperlinks) {
    if ($hl.Range.Address() -eq '$C$3') {
        $hl.TextToDisplay = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
    }
}
